$p = $ppt.ActivePresentation

# Slide 1 title: "Header" + " " + "with" + " " + "inline code" (Courier)
# -> consolidate the first four plain runs into a single run "Header with "
$s1 = $p.Slides.Item(1)
$tr1 = $s1.Shapes.Item(1).TextFrame.TextRange
$tr1.Characters(1, 12).Text = "Header with "

# Slide 2 title: "Syntax" + " " + "highlighting"
# -> consolidate into a single run "Syntax highlighting"
$s2 = $p.Slides.Item(2)
$tr2 = $s2.Shapes.Item(1).TextFrame.TextRange
$tr2.Characters(1, $tr2.Text.Length).Text = "Syntax highlighting"

# Slide 3 title: "Two" + " " + "column" + " " + "slide"
# -> consolidate into a single run "Two column slide"
$s3 = $p.Slides.Item(3)
$tr3 = $s3.Shapes.Item(1).TextFrame.TextRange
$tr3.Characters(1, $tr3.Text.Length).Text = "Two column slide"
